# Update TPM-derived NATMI ligand-receptor metrics for Efnb1-Erbb2
# (commit: "update scripts wuth new tpm")
# Applies recomputed values to columns G-T across rows 2-10 of Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 13.16594766666667
$ws.Range("H2").Value = 39.497843
$ws.Range("I2").Value = 0.6940777873489595
$ws.Range("J2").Value = 0.6940777873489595
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.451416666666667
$ws.Range("N2").Value = 7.35425
$ws.Range("O2").Value = 0.2191928499183569
$ws.Range("P2").Value = 0.2191928499183569
$ws.Range("Q2").Value = 32.27522354252778
$ws.Range("R2").Value = 290.47701188275
$ws.Range("S2").Value = 0.1521368882740457
$ws.Range("T2").Value = 0.1521368882740457
$ws.Range("G3").Value = 13.16594766666667
$ws.Range("H3").Value = 39.497843
$ws.Range("I3").Value = 0.6940777873489595
$ws.Range("J3").Value = 0.6940777873489595
$ws.Range("O3").Value = 0.4446889938320204
$ws.Range("P3").Value = 0.4446889938320204
$ws.Range("Q3").Value = 65.47858056581724
$ws.Range("R3").Value = 589.3072250923551
$ws.Range("S3").Value = 0.3086487528973638
$ws.Range("T3").Value = 0.3086487528973638
$ws.Range("G4").Value = 13.16594766666667
$ws.Range("H4").Value = 39.497843
$ws.Range("I4").Value = 0.6940777873489595
$ws.Range("J4").Value = 0.6940777873489595
$ws.Range("O4").Value = 0.3361181562496228
$ws.Range("P4").Value = 0.3361181562496228
$ws.Range("Q4").Value = 49.49198221429
$ws.Range("R4").Value = 445.42783992861
$ws.Range("S4").Value = 0.23329214617755
$ws.Range("T4").Value = 0.23329214617755
$ws.Range("I5").Value = 0.1706596770095176
$ws.Range("J5").Value = 0.1706596770095176
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 2.451416666666667
$ws.Range("N5").Value = 7.35425
$ws.Range("O5").Value = 0.2191928499183569
$ws.Range("P5").Value = 0.2191928499183569
$ws.Range("Q5").Value = 7.935824090000001
$ws.Range("R5").Value = 71.42241681
$ws.Range("S5").Value = 0.03740738096986246
$ws.Range("T5").Value = 0.03740738096986246
$ws.Range("I6").Value = 0.1706596770095176
$ws.Range("J6").Value = 0.1706596770095176
$ws.Range("O6").Value = 0.4446889938320204
$ws.Range("P6").Value = 0.4446889938320204
$ws.Range("S6").Value = 0.07589048005705996
$ws.Range("T6").Value = 0.07589048005705996
$ws.Range("I7").Value = 0.1706596770095176
$ws.Range("J7").Value = 0.1706596770095176
$ws.Range("O7").Value = 0.3361181562496228
$ws.Range("P7").Value = 0.3361181562496228
$ws.Range("S7").Value = 0.0573618159825952
$ws.Range("T7").Value = 0.0573618159825952
$ws.Range("H8").Value = 7.697376999999999
$ws.Range("I8").Value = 0.1352625356415228
$ws.Range("J8").Value = 0.1352625356415228
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 2.451416666666667
$ws.Range("N8").Value = 7.35425
$ws.Range("O8").Value = 0.2191928499183569
$ws.Range("P8").Value = 0.2191928499183569
$ws.Range("Q8").Value = 6.289826089138889
$ws.Range("R8").Value = 56.60843480225
$ws.Range("S8").Value = 0.0296485806744487
$ws.Range("T8").Value = 0.02964858067444871
$ws.Range("H9").Value = 7.697376999999999
$ws.Range("I9").Value = 0.1352625356415228
$ws.Range("J9").Value = 0.1352625356415228
$ws.Range("O9").Value = 0.4446889938320204
$ws.Range("P9").Value = 0.4446889938320204
$ws.Range("S9").Value = 0.06014976087759655
$ws.Range("T9").Value = 0.06014976087759657
$ws.Range("H10").Value = 7.697376999999999
$ws.Range("I10").Value = 0.1352625356415228
$ws.Range("J10").Value = 0.1352625356415228
$ws.Range("O10").Value = 0.3361181562496228
$ws.Range("P10").Value = 0.3361181562496228
$ws.Range("Q10").Value = 9.64504430231
$ws.Range("S10").Value = 0.04546419408947752
$ws.Range("T10").Value = 0.04546419408947753